$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.937.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.058.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.055.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("E11").Value = "  -11.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.15%  "

$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.560.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.986.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.066.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.685"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  -1.25%  "

$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("E34").Value = "  +1.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0412"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "451.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0814"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.019.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.269"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("E46").Value = "  +8.62%  "

$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("E48").Value = "  +1.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
